# Apply updated NATMI edge-weight statistics (Tgfb1-Itgav, YoungD4) per Dr Hou advice
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ "E"=3; "G"=48.746633; "H"=146.239899; "I"=0.1549390820513319; "J"=0.1549390820513319; "K"=3; "M"=18.382477; "N"=55.147431; "O"=0.06380158579420245; "P"=0.06380158579420243; "Q"=896.0838599499408; "R"=8064.754739549468; "S"=0.009885359136373027; "T"=0.009885359136373026 }
    3 = @{ "E"=3; "G"=48.746633; "H"=146.239899; "I"=0.1549390820513319; "J"=0.1549390820513319; "K"=3; "M"=82.64333833333332; "N"=247.930015; "O"=0.2868370808239535; "P"=0.2868370808239535; "Q"=4028.584483629831; "R"=36257.26035266848; "S"=0.04444227400114706; "T"=0.04444227400114706 }
    4 = @{ "E"=3; "G"=48.746633; "H"=146.239899; "I"=0.1549390820513319; "J"=0.1549390820513319; "K"=3; "M"=91.25099666666667; "N"=273.75299; "O"=0.3167123936907314; "P"=0.3167123936907314; "Q"=4448.178845394223; "R"=40033.60960854801; "S"=0.04907112755272198; "T"=0.04907112755272198 }
    5 = @{ "E"=3; "G"=48.746633; "H"=146.239899; "I"=0.1549390820513319; "J"=0.1549390820513319; "K"=3; "M"=81.28845566666666; "N"=243.865367; "O"=0.2821345773094157; "P"=0.2821345773094157; "Q"=3962.53851551977; "R"=35662.84663967793; "S"=0.04371367242326141; "T"=0.04371367242326141 }
    6 = @{ "E"=3; "G"=48.746633; "H"=146.239899; "I"=0.1549390820513319; "J"=0.1549390820513319; "K"=3; "M"=14.55416966666667; "N"=43.662509; "O"=0.0505143623816971; "P"=0.0505143623816971; "Q"=709.4667673607323; "R"=6385.20090624659; "S"=0.007826648937828484; "T"=0.007826648937828482 }
    7 = @{ "E"=3; "G"=18.65566266666667; "H"=55.966988; "I"=0.05929622356958762; "J"=0.05929622356958761; "K"=3; "M"=18.382477; "N"=55.147431; "O"=0.06380158579420245; "P"=0.06380158579420243; "Q"=342.9372898897586; "R"=3086.435609007828; "S"=0.003783193095347254; "T"=0.003783193095347252 }
    8 = @{ "E"=3; "G"=18.65566266666667; "H"=55.966988; "I"=0.05929622356958762; "J"=0.05929622356958761; "K"=3; "M"=82.64333833333332; "N"=247.930015; "O"=0.2868370808239535; "P"=0.2868370808239535; "Q"=1541.766241593869; "R"=13875.89617434482; "S"=0.01700835567258502; "T"=0.01700835567258501 }
    9 = @{ "E"=3; "G"=18.65566266666667; "H"=55.966988; "I"=0.05929622356958762; "J"=0.05929622356958761; "K"=3; "M"=91.25099666666667; "N"=273.75299; "O"=0.3167123936907314; "P"=0.3167123936907314; "Q"=1702.347811810458; "R"=15321.13030629412; "S"=0.01877984890354486; "T"=0.01877984890354486 }
    10 = @{ "E"=3; "G"=18.65566266666667; "H"=55.966988; "I"=0.05929622356958762; "J"=0.05929622356958761; "K"=3; "M"=81.28845566666666; "N"=243.865367; "O"=0.2821345773094157; "P"=0.2821345773094157; "Q"=1516.490007611622; "R"=13648.4100685046; "S"=0.01672951497285021; "T"=0.01672951497285021 }
    11 = @{ "E"=3; "G"=18.65566266666667; "H"=55.966988; "I"=0.05929622356958762; "J"=0.05929622356958761; "K"=3; "M"=14.55416966666667; "N"=43.662509; "O"=0.0505143623816971; "P"=0.0505143623816971; "Q"=271.5176796947658; "R"=2443.659117252892; "S"=0.002995310925260278; "T"=0.002995310925260277 }
    12 = @{ "E"=3; "G"=97.66137466666667; "H"=292.984124; "I"=0.310412490288807; "J"=0.310412490288807; "K"=3; "M"=18.382477; "N"=55.147431; "O"=0.06380158579420245; "P"=0.06380158579420243; "Q"=1795.257973598383; "R"=16157.32176238544; "S"=0.01980480913075336; "T"=0.01980480913075335 }
    13 = @{ "E"=3; "G"=97.66137466666667; "H"=292.984124; "I"=0.310412490288807; "J"=0.310412490288807; "K"=3; "M"=82.64333833333332; "N"=247.930015; "O"=0.2868370808239535; "P"=0.2868370808239535; "Q"=8071.062028675761; "R"=72639.55825808186; "S"=0.08903781256573522; "T"=0.08903781256573522 }
    14 = @{ "E"=3; "G"=97.66137466666667; "H"=292.984124; "I"=0.310412490288807; "J"=0.310412490288807; "K"=3; "M"=91.25099666666667; "N"=273.75299; "O"=0.3167123936907314; "P"=0.3167123936907314; "Q"=8911.697774170085; "R"=80205.27996753076; "S"=0.09831148283086899; "T"=0.09831148283086899 }
    15 = @{ "E"=3; "G"=97.66137466666667; "H"=292.984124; "I"=0.310412490288807; "J"=0.310412490288807; "K"=3; "M"=81.28845566666666; "N"=243.865367; "O"=0.2821345773094157; "P"=0.2821345773094157; "Q"=7938.742324937057; "R"=71448.6809244335; "S"=0.08757809673919566; "T"=0.08757809673919566 }
    16 = @{ "E"=3; "G"=97.66137466666667; "H"=292.984124; "I"=0.310412490288807; "J"=0.310412490288807; "K"=3; "M"=14.55416966666667; "N"=43.662509; "O"=0.0505143623816971; "P"=0.0505143623816971; "Q"=1421.380216778569; "R"=12792.42195100712; "S"=0.01568028902225383; "T"=0.01568028902225383 }
    17 = @{ "E"=3; "G"=138.7199146666667; "H"=416.159744; "I"=0.4409152985128724; "J"=0.4409152985128724; "K"=3; "M"=18.382477; "N"=55.147431; "O"=0.06380158579420245; "P"=0.06380158579420243; "Q"=2550.015640801962; "R"=22950.14076721766; "S"=0.02813109524604541; "T"=0.02813109524604541 }
    18 = @{ "E"=3; "G"=138.7199146666667; "H"=416.159744; "I"=0.4409152985128724; "J"=0.4409152985128724; "K"=3; "M"=82.64333833333332; "N"=247.930015; "O"=0.2868370808239535; "P"=0.2868370808239535; "Q"=11464.27684136846; "R"=103178.4915723162; "S"=0.1264708571160544; "T"=0.1264708571160544 }
    19 = @{ "E"=3; "G"=138.7199146666667; "H"=416.159744; "I"=0.4409152985128724; "J"=0.4409152985128724; "K"=3; "M"=91.25099666666667; "N"=273.75299; "O"=0.3167123936907314; "P"=0.3167123936907314; "Q"=12658.33047084829; "R"=113924.9742376346; "S"=0.1396433396068752; "T"=0.1396433396068752 }
    20 = @{ "E"=3; "G"=138.7199146666667; "H"=416.159744; "I"=0.4409152985128724; "J"=0.4409152985128724; "K"=3; "M"=81.28845566666666; "N"=243.865367; "O"=0.2821345773094157; "P"=0.2821345773094157; "Q"=11276.32763346512; "R"=101486.9487011861; "S"=0.1243974513751841; "T"=0.1243974513751841 }
    21 = @{ "E"=3; "G"=138.7199146666667; "H"=416.159744; "I"=0.4409152985128724; "J"=0.4409152985128724; "K"=3; "M"=14.55416966666667; "N"=43.662509; "O"=0.0505143623816971; "P"=0.0505143623816971; "Q"=2018.953174204189; "R"=18170.5785678377; "S"=0.02227255516871339; "T"=0.02227255516871339 }
    22 = @{ "E"=3; "G"=10.83447233333333; "H"=32.503417; "I"=0.03443690557740099; "J"=0.03443690557740099; "K"=3; "M"=18.382477; "N"=55.147431; "O"=0.06380158579420245; "P"=0.06380158579420243; "Q"=199.1644384746363; "R"=1792.479946271727; "S"=0.002197129185683398; "T"=0.002197129185683398 }
    23 = @{ "E"=3; "G"=10.83447233333333; "H"=32.503417; "I"=0.03443690557740099; "J"=0.03443690557740099; "K"=3; "M"=82.64333833333332; "N"=247.930015; "O"=0.2868370808239535; "P"=0.2868370808239535; "Q"=895.3969627068059; "R"=8058.572664361253; "S"=0.009877781468431823; "T"=0.009877781468431823 }
    24 = @{ "E"=3; "G"=10.83447233333333; "H"=32.503417; "I"=0.03443690557740099; "J"=0.03443690557740099; "K"=3; "M"=91.25099666666667; "N"=273.75299; "O"=0.3167123936907314; "P"=0.3167123936907314; "Q"=988.6563987740922; "R"=8897.90758896683; "S"=0.01090659479672037; "T"=0.01090659479672037 }
    25 = @{ "E"=3; "G"=10.83447233333333; "H"=32.503417; "I"=0.03443690557740099; "J"=0.03443690557740099; "K"=3; "M"=81.28845566666666; "N"=243.865367; "O"=0.2821345773094157; "P"=0.2821345773094157; "Q"=880.7175239398931; "R"=7926.457715459039; "S"=0.009715841798924287; "T"=0.009715841798924287 }
    26 = @{ "E"=3; "G"=10.83447233333333; "H"=32.503417; "I"=0.03443690557740099; "J"=0.03443690557740099; "K"=3; "M"=14.55416966666667; "N"=43.662509; "O"=0.0505143623816971; "P"=0.0505143623816971; "Q"=157.6867485881392; "R"=1419.180737293253; "S"=0.00173955832764112; "T"=0.00173955832764112 }
}

foreach ($rowNum in $updates.Keys) {
    $rowData = $updates[$rowNum]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$rowNum").Value = $rowData[$col]
    }
}
